$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 774, shifting existing rows 774:815 down to 775:816
$ws.Rows.Item(774).Insert()

# Populate the newly inserted row with the new data point.
# The date column stores plain text (not a real Excel date), so force
# text interpretation with a leading apostrophe to avoid Excel's
# automatic date-serial conversion.
$ws.Cells.Item(774, 1).Value = "'2026/02/05"
$ws.Cells.Item(774, 2).Value = "木"
$ws.Cells.Item(774, 3).Value = 10
$ws.Cells.Item(774, 4).Value = 201
